$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the monthly batch codes for existing members (June -> July 2020) ---
$ws.Range("C2").Value = "km_072020"
$ws.Range("C3").Value = "so_072020"

# --- Replace the 4th member (Abagail) with Samira Shaikh ---
$ws.Range("B4").Value = "Samira Shaikh"
$ws.Range("C4").Value = "ss_072020"
$ws.Range("D4").Value = "samirashaikh@uncc.edu"
$ws.Range("E4").Value = "Red"

# --- Remove the remaining rows (Nana, Iris, Josh) ---
$ws.Range("B5:E7").ClearContents()

# --- Rebuild the mailto hyperlinks for the email column (D2:D4) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:sroychou@uncc.edu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:samirashaikh@uncc.edu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:kmahaja2@uncc.edu") | Out-Null

# Re-apply the Hyperlink cell style so D2:D7 keep sharing the original
# hyperlink style record instead of Hyperlinks.Add minting a duplicate one.
$ws.Range("D2:D7").Style = "Hyperlink"

# --- Column width adjustments: drop the custom width on column C (back to the
# sheet's default of 14.5), resize/bestfit column D ---
$ws.Columns.Item(3).ColumnWidth = 13.67
$ws.Columns.Item(4).ColumnWidth = 19.33

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("E8").Select() | Out-Null
